$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1689184.6
$ws.Range("J17").Value = 3301386.8
$ws.Range("L17").Value = 9904160.399999999
$ws.Range("N17").Value = -9904496.399999999

$ws.Range("H29").Value = 568.8333
$ws.Range("J29").Value = 855
$ws.Range("L29").Value = 2565
$ws.Range("N29").Value = -3127

$ws.Range("H86").Value = 28574750
$ws.Range("I86").Value = 46157680
$ws.Range("J86").Value = 2488
$ws.Range("K86").Value = 46157680
$ws.Range("L86").Value = 2488
$ws.Range("M86").Value = -46156557
$ws.Range("N86").Value = -4734

$ws.Range("H89").Value = 28574750
$ws.Range("I89").Value = 46157680
$ws.Range("J89").Value = 2488
$ws.Range("K89").Value = 230788400
$ws.Range("L89").Value = 12440
$ws.Range("M89").Value = -230782784
$ws.Range("N89").Value = -23672

$ws.Range("H129").Value = 958.5780999999999
$ws.Range("I129").Value = 535.0769
$ws.Range("J129").Value = 1066.5294
$ws.Range("K129").Value = 1605.2307
$ws.Range("L129").Value = 3199.5882
$ws.Range("M129").Value = 3394.7693
$ws.Range("N129").Value = -13199.5882

$ws.Range("H135").Value = 907
$ws.Range("I135").Value = 406.25806
$ws.Range("J135").Value = 3124.5715
$ws.Range("K135").Value = 3656.32254
$ws.Range("L135").Value = 28121.1435
$ws.Range("M135").Value = -1121.32254
$ws.Range("N135").Value = -33191.1435

$ws.Range("H138").Value = 2474.258
$ws.Range("I138").Value = 2222.5715
$ws.Range("J138").Value = 3002.8
$ws.Range("K138").Value = 6667.7145
$ws.Range("L138").Value = 9008.400000000001
$ws.Range("M138").Value = -1527.7145
$ws.Range("N138").Value = -19288.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 977.63635
$ws.Range("I2").Value = 1044.5714
$ws.Range("J2").Value = 860.5
$ws.Range("K2").Value = 1044.5714
$ws.Range("L2").Value = 860.5
$ws.Range("M2").Value = -931.5714
$ws.Range("N2").Value = -1086.5

$ws.Range("H32").Value = 457100.94
$ws.Range("I32").Value = 520868.94
$ws.Range("K32").Value = 520868.94
$ws.Range("M32").Value = -520581.94

$ws.Range("H45").Value = 1792.4445
$ws.Range("I45").Value = 2133
$ws.Range("J45").Value = 1520
$ws.Range("K45").Value = 2133
$ws.Range("L45").Value = 1520
$ws.Range("M45").Value = -1756
$ws.Range("N45").Value = -2274

$ws.Range("H61").Value = 2166.4285
$ws.Range("I61").Value = 1805.3334
$ws.Range("J61").Value = 3490.4443
$ws.Range("K61").Value = 1805.3334
$ws.Range("L61").Value = 3490.4443
$ws.Range("M61").Value = -1593.3334
$ws.Range("N61").Value = -3914.4443

$ws.Range("H116").Value = 977.63635
$ws.Range("I116").Value = 1044.5714
$ws.Range("J116").Value = 860.5
$ws.Range("K116").Value = 1044.5714
$ws.Range("L116").Value = 860.5
$ws.Range("M116").Value = 1249.4286
$ws.Range("N116").Value = -5448.5

$ws.Range("H129").Value = 43389.668
$ws.Range("J129").Value = 43389.668
$ws.Range("L129").Value = 43389.668
$ws.Range("N129").Value = -53389.668

$ws.Range("H132").Value = 3407.44
$ws.Range("I132").Value = 2154.8333
$ws.Range("J132").Value = 6628.4287
$ws.Range("K132").Value = 6464.499899999999
$ws.Range("L132").Value = 19885.2861
$ws.Range("M132").Value = -3934.499899999999
$ws.Range("N132").Value = -24945.2861

$ws.Range("H136").Value = 2166.4285
$ws.Range("I136").Value = 1805.3334
$ws.Range("J136").Value = 3490.4443
$ws.Range("K136").Value = 5416.0002
$ws.Range("L136").Value = 10471.3329
$ws.Range("M136").Value = -2866.0002
$ws.Range("N136").Value = -15571.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 977.63635
$ws.Range("I3").Value = 1044.5714
$ws.Range("J3").Value = 860.5
$ws.Range("K3").Value = 1044.5714
$ws.Range("L3").Value = 860.5
$ws.Range("M3").Value = -930.5714
$ws.Range("N3").Value = -1088.5

$ws.Range("H58").Value = 43000
$ws.Range("J58").Value = 43000
$ws.Range("L58").Value = 43000
$ws.Range("N58").Value = -43588

$ws.Range("H94").Value = 1368.7931
$ws.Range("I94").Value = 1035.6818
$ws.Range("K94").Value = 1035.6818
$ws.Range("M94").Value = -584.6818000000001

$ws.Range("H131").Value = 45648
$ws.Range("J131").Value = 45648
$ws.Range("L131").Value = 45648
$ws.Range("N131").Value = -55728

$ws.Range("H133").Value = 42857.145
$ws.Range("J133").Value = 42857.145
$ws.Range("L133").Value = 42857.145
$ws.Range("N133").Value = -52977.145

$ws.Range("H134").Value = 3077.5217
$ws.Range("I134").Value = 2654.9285
$ws.Range("J134").Value = 3734.889
$ws.Range("K134").Value = 7964.7855
$ws.Range("L134").Value = 11204.667
$ws.Range("M134").Value = -5429.7855
$ws.Range("N134").Value = -16274.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5076.217
$ws.Range("I31").Value = 1247.0286
$ws.Range("J31").Value = 10437.08
$ws.Range("K31").Value = 1247.0286
$ws.Range("L31").Value = 10437.08
$ws.Range("M31").Value = -952.0286000000001
$ws.Range("N31").Value = -11027.08

$ws.Range("H34").Value = 5076.217
$ws.Range("I34").Value = 1247.0286
$ws.Range("J34").Value = 10437.08
$ws.Range("K34").Value = 1247.0286
$ws.Range("L34").Value = 10437.08
$ws.Range("M34").Value = -1045.0286
$ws.Range("N34").Value = -10841.08

$ws.Range("H58").Value = 1543.871
$ws.Range("I58").Value = 1271.125
$ws.Range("J58").Value = 1834.8
$ws.Range("K58").Value = 1271.125
$ws.Range("L58").Value = 1834.8
$ws.Range("M58").Value = -1068.125
$ws.Range("N58").Value = -2240.8

$ws.Range("H107").Value = 451.25
$ws.Range("I107").Value = 473.33334
$ws.Range("J107").Value = 385
$ws.Range("K107").Value = 473.33334
$ws.Range("L107").Value = 385
$ws.Range("M107").Value = 1446.66666
$ws.Range("N107").Value = -4225

$ws.Range("H132").Value = 4506209
$ws.Range("I132").Value = 1356.6538
$ws.Range("J132").Value = 15154042
$ws.Range("K132").Value = 4069.9614
$ws.Range("L132").Value = 45462126
$ws.Range("M132").Value = -1539.9614
$ws.Range("N132").Value = -45467186

$ws.Range("H134").Value = 4417.2163
$ws.Range("I134").Value = 4906.4614
$ws.Range("J134").Value = 3260.818
$ws.Range("K134").Value = 14719.3842
$ws.Range("L134").Value = 9782.454000000002
$ws.Range("M134").Value = -12184.3842
$ws.Range("N134").Value = -14852.454

$ws.Range("H136").Value = 1543.871
$ws.Range("I136").Value = 1271.125
$ws.Range("J136").Value = 1834.8
$ws.Range("K136").Value = 3813.375
$ws.Range("L136").Value = 5504.4
$ws.Range("M136").Value = -1263.375
$ws.Range("N136").Value = -10604.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1622.3529
$ws.Range("J31").Value = 1536.25
$ws.Range("L31").Value = 4608.75
$ws.Range("N31").Value = -5184.75

$ws.Range("H113").Value = 776.1429000000001
$ws.Range("I113").Value = 572.6786
$ws.Range("K113").Value = 1718.0358
$ws.Range("M113").Value = 451.9642000000001

$ws.Range("H131").Value = 1269.6522
$ws.Range("J131").Value = 1625.5333
$ws.Range("L131").Value = 4876.5999
$ws.Range("N131").Value = -14956.5999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 56557724
$ws.Range("I80").Value = 101801100
$ws.Range("J80").Value = 3496.25
$ws.Range("K80").Value = 101801100
$ws.Range("L80").Value = 3496.25
$ws.Range("M80").Value = -101800102
$ws.Range("N80").Value = -5492.25

$ws.Range("H83").Value = 56557724
$ws.Range("I83").Value = 101801100
$ws.Range("J83").Value = 3496.25
$ws.Range("K83").Value = 509005500
$ws.Range("L83").Value = 17481.25
$ws.Range("M83").Value = -509000508
$ws.Range("N83").Value = -27465.25

$ws.Range("H122").Value = 1500.6666
$ws.Range("I122").Value = 1072.2858
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3216.8574
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -766.8574000000003
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 2514.577
$ws.Range("I132").Value = 2226.9
$ws.Range("J132").Value = 3473.5
$ws.Range("K132").Value = 6680.700000000001
$ws.Range("L132").Value = 10420.5
$ws.Range("M132").Value = -4150.700000000001
$ws.Range("N132").Value = -15480.5

$ws.Range("H137").Value = 28088.422
$ws.Range("J137").Value = 29371.111
$ws.Range("L137").Value = 29371.111
$ws.Range("N137").Value = -39571.111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4784.911
$ws.Range("I22").Value = 2371.4285
$ws.Range("J22").Value = 8760.058999999999
$ws.Range("K22").Value = 2371.4285
$ws.Range("L22").Value = 8760.058999999999
$ws.Range("M22").Value = -2076.4285
$ws.Range("N22").Value = -9350.058999999999

$ws.Range("H27").Value = 4784.911
$ws.Range("I27").Value = 2371.4285
$ws.Range("J27").Value = 8760.058999999999
$ws.Range("K27").Value = 2371.4285
$ws.Range("L27").Value = 8760.058999999999
$ws.Range("M27").Value = -2264.4285
$ws.Range("N27").Value = -8974.058999999999

$ws.Range("H132").Value = 1691.8125
$ws.Range("I132").Value = 1249.4445
$ws.Range("J132").Value = 3018.9167
$ws.Range("K132").Value = 3748.3335
$ws.Range("L132").Value = 9056.750100000001
$ws.Range("M132").Value = -1218.3335
$ws.Range("N132").Value = -14116.7501

$ws.Range("H133").Value = 34529.145
$ws.Range("J133").Value = 34529.145
$ws.Range("L133").Value = 34529.145
$ws.Range("N133").Value = -39589.145

$ws.Range("H136").Value = 4274592.5
$ws.Range("I136").Value = 920
$ws.Range("K136").Value = 2760
$ws.Range("M136").Value = -210

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

$ws.Range("H136").Value = 1706.614
$ws.Range("I136").Value = 1253.5834
$ws.Range("K136").Value = 3760.7502
$ws.Range("M136").Value = -1210.7502
